# Update cryptocurrency price/volume figures per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.143.74'
$ws.Range('E2').Value = '  +5.75%  '
$ws.Range('D3').Value = '1.917.26'
$ws.Range('E3').Value = '  +2.50%  '
$ws.Range('E4').Value = '  -0.80%  '
$ws.Range('D5').Value = "'329.96"
$ws.Range('E5').Value = '  +4.61%  '
$ws.Range('D7').Value = "'0.5238"
$ws.Range('E7').Value = '  +3.13%  '
$ws.Range('D8').Value = "'0.4079"
$ws.Range('E8').Value = '  +4.54%  '
$ws.Range('D9').Value = "'0.08521"
$ws.Range('E9').Value = '  +1.88%  '
$ws.Range('D10').Value = "'42.89"
$ws.Range('E10').Value = '  +1.26%  '
$ws.Range('E11').Value = '  +1.56%  '
$ws.Range('D12').Value = "'22.45"
$ws.Range('E12').Value = '  +10.25%  '
$ws.Range('D13').Value = "'6.449"
$ws.Range('E13').Value = '  +4.02%  '
$ws.Range('D14').Value = '1.918.71'
$ws.Range('E14').Value = '  +2.54%  '
$ws.Range('D15').Value = "'7.387"
$ws.Range('E15').Value = '  +1.89%  '
$ws.Range('E16').Value = '  -0.82%  '
$ws.Range('D17').Value = "'94.94"
$ws.Range('E17').Value = '  +3.99%  '
$ws.Range('E18').Value = '  +1.14%  '
$ws.Range('D19').Value = "'0.06688"
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('D20').Value = "'18.39"
$ws.Range('E20').Value = '  +4.15%  '
$ws.Range('D21').Value = "'1.000"
$ws.Range('E21').Value = '  -0.71%  '
$ws.Range('D22').Value = "'6.014"
$ws.Range('E22').Value = '  +1.76%  '
$ws.Range('D23').Value = '30.209.75'
$ws.Range('E23').Value = '  +5.77%  '
$ws.Range('D24').Value = "'11.34"
$ws.Range('E24').Value = '  +2.23%  '
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('D26').Value = '2.135.29'
$ws.Range('D27').Value = "'160.62"
$ws.Range('E27').Value = '  +2.39%  '
$ws.Range('D28').Value = "'21.14"
$ws.Range('E28').Value = '  +2.75%  '
$ws.Range('D29').Value = "'2.419"
$ws.Range('E29').Value = '  -0.31%  '
$ws.Range('D30').Value = "'129.01"
$ws.Range('E30').Value = '  +2.39%  '
$ws.Range('D31').Value = "'1.079"
$ws.Range('E31').Value = '  +3.86%  '
$ws.Range('D32').Value = "'0.1065"
$ws.Range('E32').Value = '  +2.51%  '
$ws.Range('D33').Value = "'5.997"
$ws.Range('E33').Value = '  +4.30%  '
$ws.Range('D34').Value = "'3.639"
$ws.Range('E34').Value = '  +0.51%  '
$ws.Range('D35').Value = "'0.02485"
$ws.Range('E35').Value = '  +1.33%  '
$ws.Range('D36').Value = "'0.06603"
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D37').Value = "'0.2208"
$ws.Range('E37').Value = '  +2.10%  '
$ws.Range('E38').Value = '  +2.66%  '
$ws.Range('D39').Value = "'1.227"
$ws.Range('E39').Value = '  +3.88%  '
$ws.Range('D40').Value = "'8.875"
$ws.Range('E40').Value = '  -1.19%  '
$ws.Range('D41').Value = "'0.6535"
$ws.Range('E41').Value = '  +2.65%  '
$ws.Range('E42').Value = '  +4.75%  '
$ws.Range('D43').Value = "'1.242"
$ws.Range('E43').Value = '  +0.63%  '
$ws.Range('D44').Value = "'0.6143"
$ws.Range('E44').Value = '  +2.32%  '
$ws.Range('D45').Value = "'13.21"
$ws.Range('E45').Value = '  +1.57%  '
$ws.Range('D46').Value = "'3.749"
$ws.Range('E46').Value = '  +1.85%  '
$ws.Range('D47').Value = "'2.081"
$ws.Range('E47').Value = '  +4.01%  '
$ws.Range('D48').Value = "'1.244"
$ws.Range('E48').Value = '  +2.63%  '
$ws.Range('D49').Value = "'124.46"
$ws.Range('E49').Value = '  +1.67%  '
$ws.Range('E50').Value = '  +3.99%  '
$ws.Range('D51').Value = "'79.63"
$ws.Range('E51').Value = '  +4.41%  '
